$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.507882
$ws.Cells.Item(2, 8).Value = 16.523646
$ws.Cells.Item(2, 9).Value = 0.03518866199235487
$ws.Cells.Item(2, 10).Value = 0.03518866199235487
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.048438999999999
$ws.Cells.Item(2, 14).Value = 12.145317
$ws.Cells.Item(2, 15).Value = 0.1703267688113503
$ws.Cells.Item(2, 16).Value = 0.1703267688113503
$ws.Cells.Item(2, 17).Value = 22.29832429619799
$ws.Cells.Item(2, 18).Value = 200.684918665782
$ws.Cells.Item(2, 19).Value = 0.005993571095952575
$ws.Cells.Item(2, 20).Value = 0.005993571095952576

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.507882
$ws.Cells.Item(3, 8).Value = 16.523646
$ws.Cells.Item(3, 9).Value = 0.03518866199235487
$ws.Cells.Item(3, 10).Value = 0.03518866199235487
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 11.27122833333333
$ws.Cells.Item(3, 14).Value = 33.813685
$ws.Cells.Item(3, 15).Value = 0.4742054659960562
$ws.Cells.Item(3, 16).Value = 0.4742054659960562
$ws.Cells.Item(3, 17).Value = 62.08059565505666
$ws.Cells.Item(3, 18).Value = 558.72536089551
$ws.Cells.Item(3, 19).Value = 0.01668665585786235
$ws.Cells.Item(3, 20).Value = 0.01668665585786235

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.507882
$ws.Cells.Item(4, 8).Value = 16.523646
$ws.Cells.Item(4, 9).Value = 0.03518866199235487
$ws.Cells.Item(4, 10).Value = 0.03518866199235487
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.219226333333333
$ws.Cells.Item(4, 14).Value = 9.657679
$ws.Cells.Item(4, 15).Value = 0.1354399607920677
$ws.Cells.Item(4, 16).Value = 0.1354399607920676
$ws.Cells.Item(4, 17).Value = 17.73111877529266
$ws.Cells.Item(4, 18).Value = 159.580068977634
$ws.Cells.Item(4, 19).Value = 0.004765951000569865
$ws.Cells.Item(4, 20).Value = 0.004765951000569864

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.507882
$ws.Cells.Item(5, 8).Value = 16.523646
$ws.Cells.Item(5, 9).Value = 0.03518866199235487
$ws.Cells.Item(5, 10).Value = 0.03518866199235487
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.235549333333334
$ws.Cells.Item(5, 14).Value = 12.706648
$ws.Cells.Item(5, 15).Value = 0.1781989137264352
$ws.Cells.Item(5, 16).Value = 0.1781989137264352
$ws.Cells.Item(5, 17).Value = 23.32890593317867
$ws.Cells.Item(5, 18).Value = 209.960153398608
$ws.Cells.Item(5, 19).Value = 0.006270581342524336
$ws.Cells.Item(5, 20).Value = 0.006270581342524334

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.507882
$ws.Cells.Item(6, 8).Value = 16.523646
$ws.Cells.Item(6, 9).Value = 0.03518866199235487
$ws.Cells.Item(6, 10).Value = 0.03518866199235487
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9942166666666665
$ws.Cells.Item(6, 14).Value = 2.98265
$ws.Cells.Item(6, 15).Value = 0.04182889067409059
$ws.Cells.Item(6, 16).Value = 0.04182889067409059
$ws.Cells.Item(6, 17).Value = 5.476028082433333
$ws.Cells.Item(6, 18).Value = 49.28425274189999
$ws.Cells.Item(6, 19).Value = 0.001471902695445739
$ws.Cells.Item(6, 20).Value = 0.001471902695445739

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 149.9875183333334
$ws.Cells.Item(7, 8).Value = 449.9625550000001
$ws.Cells.Item(7, 9).Value = 0.9582376829612175
$ws.Cells.Item(7, 10).Value = 0.9582376829612176
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.048438999999999
$ws.Cells.Item(7, 14).Value = 12.145317
$ws.Cells.Item(7, 15).Value = 0.1703267688113503
$ws.Cells.Item(7, 16).Value = 0.1703267688113503
$ws.Cells.Item(7, 17).Value = 607.2153187338816
$ws.Cells.Item(7, 18).Value = 5464.937868604935
$ws.Cells.Item(7, 19).Value = 0.1632135282920592
$ws.Cells.Item(7, 20).Value = 0.1632135282920593

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 149.9875183333334
$ws.Cells.Item(8, 8).Value = 449.9625550000001
$ws.Cells.Item(8, 9).Value = 0.9582376829612175
$ws.Cells.Item(8, 10).Value = 0.9582376829612176
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 11.27122833333333
$ws.Cells.Item(8, 14).Value = 33.813685
$ws.Cells.Item(8, 15).Value = 0.4742054659960562
$ws.Cells.Item(8, 16).Value = 0.4742054659960562
$ws.Cells.Item(8, 17).Value = 1690.54356628502
$ws.Cells.Item(8, 18).Value = 15214.89209656518
$ws.Cells.Item(8, 19).Value = 0.4544015469836053
$ws.Cells.Item(8, 20).Value = 0.4544015469836053

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 149.9875183333334
$ws.Cells.Item(9, 8).Value = 449.9625550000001
$ws.Cells.Item(9, 9).Value = 0.9582376829612175
$ws.Cells.Item(9, 10).Value = 0.9582376829612176
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.219226333333333
$ws.Cells.Item(9, 14).Value = 9.657679
$ws.Cells.Item(9, 15).Value = 0.1354399607920677
$ws.Cells.Item(9, 16).Value = 0.1354399607920676
$ws.Cells.Item(9, 17).Value = 482.8437686899828
$ws.Cells.Item(9, 18).Value = 4345.593918209845
$ws.Cells.Item(9, 19).Value = 0.1297836742097491
$ws.Cells.Item(9, 20).Value = 0.129783674209749

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 149.9875183333334
$ws.Cells.Item(10, 8).Value = 449.9625550000001
$ws.Cells.Item(10, 9).Value = 0.9582376829612175
$ws.Cells.Item(10, 10).Value = 0.9582376829612176
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.235549333333334
$ws.Cells.Item(10, 14).Value = 12.706648
$ws.Cells.Item(10, 15).Value = 0.1781989137264352
$ws.Cells.Item(10, 16).Value = 0.1781989137264352
$ws.Cells.Item(10, 17).Value = 635.2795332850714
$ws.Cells.Item(10, 18).Value = 5717.515799565642
$ws.Cells.Item(10, 19).Value = 0.1707569141954252
$ws.Cells.Item(10, 20).Value = 0.1707569141954252

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 149.9875183333334
$ws.Cells.Item(11, 8).Value = 449.9625550000001
$ws.Cells.Item(11, 9).Value = 0.9582376829612175
$ws.Cells.Item(11, 10).Value = 0.9582376829612176
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.9942166666666665
$ws.Cells.Item(11, 14).Value = 2.98265
$ws.Cells.Item(11, 15).Value = 0.04182889067409059
$ws.Cells.Item(11, 16).Value = 0.04182889067409059
$ws.Cells.Item(11, 17).Value = 149.1200905189722
$ws.Cells.Item(11, 18).Value = 1342.08081467075
$ws.Cells.Item(11, 19).Value = 0.04008201928037865
$ws.Cells.Item(11, 20).Value = 0.04008201928037865

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.028937
$ws.Cells.Item(12, 8).Value = 3.086811
$ws.Cells.Item(12, 9).Value = 0.006573655046427582
$ws.Cells.Item(12, 10).Value = 0.006573655046427582
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.048438999999999
$ws.Cells.Item(12, 14).Value = 12.145317
$ws.Cells.Item(12, 15).Value = 0.1703267688113503
$ws.Cells.Item(12, 16).Value = 0.1703267688113503
$ws.Cells.Item(12, 17).Value = 4.165588679342999
$ws.Cells.Item(12, 18).Value = 37.490298114087
$ws.Cells.Item(12, 19).Value = 0.001119669423338437
$ws.Cells.Item(12, 20).Value = 0.001119669423338437

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.028937
$ws.Cells.Item(13, 8).Value = 3.086811
$ws.Cells.Item(13, 9).Value = 0.006573655046427582
$ws.Cells.Item(13, 10).Value = 0.006573655046427582
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 11.27122833333333
$ws.Cells.Item(13, 14).Value = 33.813685
$ws.Cells.Item(13, 15).Value = 0.4742054659960562
$ws.Cells.Item(13, 16).Value = 0.4742054659960562
$ws.Cells.Item(13, 17).Value = 11.597383867615
$ws.Cells.Item(13, 18).Value = 104.376454808535
$ws.Cells.Item(13, 19).Value = 0.003117263154588518
$ws.Cells.Item(13, 20).Value = 0.003117263154588518

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.028937
$ws.Cells.Item(14, 8).Value = 3.086811
$ws.Cells.Item(14, 9).Value = 0.006573655046427582
$ws.Cells.Item(14, 10).Value = 0.006573655046427582
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.219226333333333
$ws.Cells.Item(14, 14).Value = 9.657679
$ws.Cells.Item(14, 15).Value = 0.1354399607920677
$ws.Cells.Item(14, 16).Value = 0.1354399607920676
$ws.Cells.Item(14, 17).Value = 3.312381085741
$ws.Cells.Item(14, 18).Value = 29.811429771669
$ws.Cells.Item(14, 19).Value = 0.0008903355817487294
$ws.Cells.Item(14, 20).Value = 0.0008903355817487292

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.028937
$ws.Cells.Item(15, 8).Value = 3.086811
$ws.Cells.Item(15, 9).Value = 0.006573655046427582
$ws.Cells.Item(15, 10).Value = 0.006573655046427582
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 4.235549333333334
$ws.Cells.Item(15, 14).Value = 12.706648
$ws.Cells.Item(15, 15).Value = 0.1781989137264352
$ws.Cells.Item(15, 16).Value = 0.1781989137264352
$ws.Cells.Item(15, 17).Value = 4.358113424392001
$ws.Cells.Item(15, 18).Value = 39.223020819528
$ws.Cells.Item(15, 19).Value = 0.001171418188485694
$ws.Cells.Item(15, 20).Value = 0.001171418188485694

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.028937
$ws.Cells.Item(16, 8).Value = 3.086811
$ws.Cells.Item(16, 9).Value = 0.006573655046427582
$ws.Cells.Item(16, 10).Value = 0.006573655046427582
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.9942166666666665
$ws.Cells.Item(16, 14).Value = 2.98265
$ws.Cells.Item(16, 15).Value = 0.04182889067409059
$ws.Cells.Item(16, 16).Value = 0.04182889067409059
$ws.Cells.Item(16, 17).Value = 1.02298631435
$ws.Cells.Item(16, 18).Value = 9.206876829149998
$ws.Cells.Item(16, 19).Value = 0.0002749686982662032
$ws.Cells.Item(16, 20).Value = 0.0002749686982662032

